$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 315, shifting existing rows 315-329 down to 317-331.
$ws.Range("A315:R316").EntireRow.Insert()

# New row 315
$ws.Range("A315").Value = 10
$ws.Range("B315").Value = "Vega Modelo de Temuco"
$ws.Range("C315").Value = "La Araucanía"
$ws.Range("D315").Value = 44753
$ws.Range("E315").Value = 9
$ws.Range("F315").Value = 100112044
$ws.Range("G315").Value = "Perejil"
$ws.Range("H315").Value = "Sin especificar"
$ws.Range("I315").Value = "Primera"
$ws.Range("J315").Value = 30
$ws.Range("K315").Value = 5000
$ws.Range("L315").Value = 5000
$ws.Range("M315").Value = 5000
$ws.Range("N315").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O315").Value = "Provincia de Cautín"
$ws.Range("P315").Value = 1667
$ws.Range("Q315").Value = 3
$ws.Range("R315").Value = "Hortaliza"

# New row 316
$ws.Range("A316").Value = 10
$ws.Range("B316").Value = "Vega Modelo de Temuco"
$ws.Range("C316").Value = "La Araucanía"
$ws.Range("D316").Value = 44753
$ws.Range("E316").Value = 9
$ws.Range("F316").Value = 100112044
$ws.Range("G316").Value = "Perejil"
$ws.Range("H316").Value = "Sin especificar"
$ws.Range("I316").Value = "Primera"
$ws.Range("J316").Value = 50
$ws.Range("K316").Value = 4300
$ws.Range("L316").Value = 4300
$ws.Range("M316").Value = 4300
$ws.Range("N316").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O316").Value = "Región Metropolitana"
$ws.Range("P316").Value = 1433
$ws.Range("Q316").Value = 3
$ws.Range("R316").Value = "Hortaliza"
